$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update existing row 48 (A48 date changes from 45896 to 45927)
$ws.Range("A48").Value = 45927

# Add new row 49
$ws.Range("A49").Value = 45932
$ws.Range("B49").Value = 0.85416666666666663
$ws.Range("C49").Value = 48
$ws.Range("D49").Value = "Flamengo"
$ws.Range("E49").Value = "Brasileiro"
$ws.Range("F49").Value = "Fora"
$ws.Range("G49").Value = "Maracanã"
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = "Empate"

# Add new row 50
$ws.Range("A50").Value = 45935
$ws.Range("B50").Value = 0.85416666666666663
$ws.Range("C50").Value = 49
$ws.Range("D50").Value = "Sport"
$ws.Range("E50").Value = "Brasileiro"
$ws.Range("F50").Value = "Casa"
$ws.Range("G50").Value = "Mineirão"
$ws.Range("H50").Value = 1
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = "Empate"

# Apply the same number formats as the row above for the date/time columns
$ws.Range("A49:A50").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B49:B50").NumberFormat = $ws.Range("B48").NumberFormat

# Update sheet view to match the new scroll position / selection
$ws.Application.ActiveWindow.ScrollRow = 29
$sel = $ws.Range("J51")
$sel.Select()
